# Add the newly-found event points (rows 27-48) that Michelle found.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2001-08-12T08:00:00.000", "2001-08-12T20:00:00.000"),
    @("2001-08-15T00:00:00.000", "2001-08-15T09:00:00.000"),
    @("2001-08-17T04:00:00.000", "2001-08-17T16:00:00.000"),
    @("2001-08-19T14:00:00.000", "2001-08-19T20:30:00.000"),
    @("2001-08-22T04:00:00.000", "2001-08-22T10:15:00.000"),
    @("2001-08-24T12:00:00.000", "2001-08-24T15:30:00.000"),
    @("2001-08-26T20:00:00.000", "2001-08-27T05:45:00.000"),
    @("2001-08-29T04:00:00.000", "2001-08-29T11:00:00.000"),
    @("2001-08-31T08:00:00.000", "2001-08-31T17:45:00.000"),
    @("2001-09-05T04:00:00.000", "2001-09-05T11:15:00.000"),
    @("2001-09-10T00:00:00.000", "2001-09-10T08:30:00.000"),
    @("2001-09-12T08:00:00.000", "2001-09-12T14:30:00.000"),
    @("2001-09-17T00:00:00.000", "2001-09-17T09:45:00.000"),
    @("2001-09-19T08:00:00.000", "2001-09-19T16:45:00.000"),
    @("2001-09-21T20:00:00.000", "2001-09-22T05:15:00.000"),
    @("2001-09-24T00:00:00.000", "2001-09-24T09:00:00.000"),
    @("2001-09-29T00:00:00.000", "2001-09-29T07:00:00.000"),
    @("2001-10-01T08:00:00.000", "2001-10-01T10:30:00.000"),
    @("2001-10-03T20:00:00.000", "2001-10-04T00:45:00.000"),
    @("2001-10-06T00:00:00.000", "2001-10-06T08:00:00.000"),
    @("2001-10-08T04:00:00.000", "2001-10-08T14:00:00.000"),
    @("2001-10-11T00:00:00.000", "2001-10-11T04:30:00.000")
)

$startRow = 27
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $aVal = $data[$i][0]
    $cVal = $data[$i][1]

    $rangeA = $ws.Range("A" + $row)
    $rangeB = $ws.Range("B" + $row)
    $rangeC = $ws.Range("C" + $row)

    # Column A & C: right-aligned, wrapped text, formatted as Text ("@")
    $rangeA.NumberFormat = "@"
    $rangeA.Font.Name = "Arial"
    $rangeA.Font.Size = 11
    $rangeA.Font.Color = 0
    $rangeA.HorizontalAlignment = -4152
    $rangeA.WrapText = $true
    $rangeA.Value = $aVal

    $rangeC.NumberFormat = "@"
    $rangeC.Font.Name = "Arial"
    $rangeC.Font.Size = 11
    $rangeC.Font.Color = 0
    $rangeC.HorizontalAlignment = -4152
    $rangeC.WrapText = $true
    $rangeC.Value = $cVal

    # Column B: numeric 1, general format, default font
    $rangeB.NumberFormat = "General"
    $rangeB.Font.Name = "Arial"
    $rangeB.Font.Size = 10
    $rangeB.Font.Color = 0
    $rangeB.Value = 1

    $ws.Rows.Item($row).RowHeight = 14.15
}

# Column width adjustments
$ws.Columns.Item(1).ColumnWidth = 28.06
$ws.Columns.Item(3).ColumnWidth = 25.98

# Update selection to match the author's final cursor position
$ws.Range("C53").Select()
